$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old rows 13-24 down to 14-25)
$ws.Rows.Item(13).Insert()

# Fix column definition: column A alone should be 30.7109375 wide (previously merged with column B span)
$ws.Columns.Item(1).ColumnWidth = 29.83

# Clear cell A13 (inserted blank row separator - no label here)
$ws.Range("A13").ClearContents()

# Set final cell values
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Range("B2").Value = 'LOT2022'
$ws.Range("C2").Value = 'LOT2022'
$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Modelagem e Simulação de Processos Biotecnológicos'
$ws.Range("C3").Value = ' Modelagem e Simulação de Processos Biotecnológicos'
$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Modeling and Simulation of Biotechnological Processes'
$ws.Range("C4").Value = 'Modeling and Simulation of Biotechnological Processes'
$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '3'
$ws.Range("C5").Value = '3'
$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '0'
$ws.Range("C6").Value = '0'
$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '45 h'
$ws.Range("C7").Value = '45 h'
$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2018'
$ws.Range("C8").Value = '01/01/2018'
$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EB-8'
$ws.Range("C9").Value = 'EB-8'
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Capacitar o aluno para a análise e simulação de processos biotecnológicos através do desenvolvimento de modelos matemáticos baseados em princípios de conservação de massa, energia e quantidade de movimento, além de equações constitutivas, condições iniciais e de contorno. Serão ainda apresentadas e discutidas ferramentas computacionais aplicadas à resolução de modelos matemáticos e à simulação de processos, enfatizando o uso destas em problemas de engenharia bioquímica.'
$ws.Range("C10").Value = 'Capacitar o aluno para a análise e simulação de processos biotecnológicos através do desenvolvimento de modelos matemáticos baseados em princípios de conservação de massa, energia e quantidade de movimento, além de equações constitutivas, condições iniciais e de contorno. Serão ainda apresentadas e discutidas ferramentas computacionais aplicadas à resolução de modelos matemáticos e à simulação de processos, enfatizando o uso destas em problemas de engenharia bioquímica.'
$ws.Range("A11").Value = 'Objectives:'
$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("B13").Value = '6007846 - Júlio César dos Santos'
$ws.Range("C13").Value = '6007846 - Júlio César dos Santos'
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = 'Introdução à modelagem e simulação de bioprocessos; Estudo de problemas ligados à indústria de bioprocessos envolvendo o desenvolvimento e a resolução de modelos fenomenológicos: programas computacionais e equações algébricas; Modelagem matemática e simulação de processos fermentativos; Desenvolvimento e resolução de modelos: equações diferenciais; Ajuste de parâmetros e otimização de bioprocessos; Utilização de simuladores de processos aplicada à biotecnologia.'
$ws.Range("C14").Value = 'Introdução à modelagem e simulação de bioprocessos; Estudo de problemas ligados à indústria de bioprocessos envolvendo o desenvolvimento e a resolução de modelos fenomenológicos: programas computacionais e equações algébricas; Modelagem matemática e simulação de processos fermentativos; Desenvolvimento e resolução de modelos: equações diferenciais; Ajuste de parâmetros e otimização de bioprocessos; Utilização de simuladores de processos aplicada à biotecnologia.'
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = 'Introduction to modeling and simulation of bioprocesses; study of problems of the industry of bioprocesses related to the construction and solution of phenomenological models: computational software and algebraic equations; mathematical modeling and simulation of fermentative processes; constructing and solving models: differential equations; adjustment of parameters and bioprocesses optimization; use of process simulators applied to biotechnology.'
$ws.Range("C15").Value = 'Introduction to modeling and simulation of bioprocesses; study of problems of the industry of bioprocesses related to the construction and solution of phenomenological models: computational software and algebraic equations; mathematical modeling and simulation of fermentative processes; constructing and solving models: differential equations; adjustment of parameters and bioprocesses optimization; use of process simulators applied to biotechnology.'
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = '1. Introdução à modelagem e simulação de bioprocessos1.1. Definição de modelo matemático1.2. Conceituação de variáveis dependentes e independentes de um sistema1.3. Definição e classificação de volume de controle2. Estudo de problemas ligados à indústria de bioprocessos envolvendo o desenvolvimento e a resolução de modelos fenomenológicos: programas computacionais e equações algébricas2.1 Introdução ao programa computacional utilizado para a resolução dos modelos matemáticos2.2 Problemas envolvendo sistemas de equações lineares2.3 Problemas envolvendo equações não lineares2.4 Problemas envolvendo sistemas de equações não lineares3. Modelagem matemática e simulação de processos fermentativos3.1. Objetivos3.2. Diferenças entre processos químicos e fermentativos3.3. Interações entre a população microbiana e o meio de cultura3.4. Formulação e classificação de modelos matemáticos de processos fermentativos 3.5. Modelos cinéticos de crescimento celular, consumo de substrato e formação de produtos em processos fermentativos3.6. Modelagem de processo fermentativo em reator: descontínuo, contínuo, contínuo com reciclo de células, descontínuo alimentado e tubular.4. Desenvolvimento e resolução de modelos: equações diferenciais5. Ajuste de parâmetros e otimização de bioprocessos6. Utilização de simuladores de processos aplicada à biotecnologia6.1. Projetos auxiliados por pacotes computacionais de simulação de projetos6.2. Classificação dos pacotes computacionais de simulação de processos6.3. Síntese e análise de processos6.4. Desenvolvimento de fluxogramas de simulação: conceitos e limitações, convergência.6.5. Exemplos de aplicação'
$ws.Range("C16").Value = '1. Introdução à modelagem e simulação de bioprocessos1.1. Definição de modelo matemático1.2. Conceituação de variáveis dependentes e independentes de um sistema1.3. Definição e classificação de volume de controle2. Estudo de problemas ligados à indústria de bioprocessos envolvendo o desenvolvimento e a resolução de modelos fenomenológicos: programas computacionais e equações algébricas2.1 Introdução ao programa computacional utilizado para a resolução dos modelos matemáticos2.2 Problemas envolvendo sistemas de equações lineares2.3 Problemas envolvendo equações não lineares2.4 Problemas envolvendo sistemas de equações não lineares3. Modelagem matemática e simulação de processos fermentativos3.1. Objetivos3.2. Diferenças entre processos químicos e fermentativos3.3. Interações entre a população microbiana e o meio de cultura3.4. Formulação e classificação de modelos matemáticos de processos fermentativos 3.5. Modelos cinéticos de crescimento celular, consumo de substrato e formação de produtos em processos fermentativos3.6. Modelagem de processo fermentativo em reator: descontínuo, contínuo, contínuo com reciclo de células, descontínuo alimentado e tubular.4. Desenvolvimento e resolução de modelos: equações diferenciais5. Ajuste de parâmetros e otimização de bioprocessos6. Utilização de simuladores de processos aplicada à biotecnologia6.1. Projetos auxiliados por pacotes computacionais de simulação de projetos6.2. Classificação dos pacotes computacionais de simulação de processos6.3. Síntese e análise de processos6.4. Desenvolvimento de fluxogramas de simulação: conceitos e limitações, convergência.6.5. Exemplos de aplicação'
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = '1. Introduction to modeling and simulation of bioprocesses.1.1. Definition of mathematical model.1.2. Concepts of dependent and independent variables of a system.1.3. Definition and classification of control volumes.2. Study of problems of the industry of bioprocesses related to the construction and solution of phenomenological models: computational software and algebraic equations.2.1. Introduction to computational software/packages used to solving mathematical models.2.2. Solving of problems using systems of linear equations.2.3. Solving of problems using non-linear equations.2.4. Solving of problems using systems of non-linear equations.3. Mathematical modeling and simulation of fermentative processes3.1. Objectives3.2. Differences between chemical and fermentative processes3.3. Interactions between the microbial population and the culture medium.3.4. Construction and classification of mathematical models for fermentative processes.3.5. Kinetic models of cellular growth, substrate consumption and formation of products in fermentative processes.3.6. Modeling of fermentative process in reactors: batch, continuous, continuous with cells recycle, fed-batch and tubular.4. Constructing and solving models: differential equations.5. Adjustment of parameters and bioprocesses optimization.6. Use of process simulators applied to biotechnology.6.1. Process design aided by simulation software6.2. Classification of process simulation software 6.3. Synthesis and analysis of process6.4. Process flowsheeting: concepts and limitation, convergence6.5. Application examples.'
$ws.Range("C17").Value = '1. Introduction to modeling and simulation of bioprocesses.1.1. Definition of mathematical model.1.2. Concepts of dependent and independent variables of a system.1.3. Definition and classification of control volumes.2. Study of problems of the industry of bioprocesses related to the construction and solution of phenomenological models: computational software and algebraic equations.2.1. Introduction to computational software/packages used to solving mathematical models.2.2. Solving of problems using systems of linear equations.2.3. Solving of problems using non-linear equations.2.4. Solving of problems using systems of non-linear equations.3. Mathematical modeling and simulation of fermentative processes3.1. Objectives3.2. Differences between chemical and fermentative processes3.3. Interactions between the microbial population and the culture medium.3.4. Construction and classification of mathematical models for fermentative processes.3.5. Kinetic models of cellular growth, substrate consumption and formation of products in fermentative processes.3.6. Modeling of fermentative process in reactors: batch, continuous, continuous with cells recycle, fed-batch and tubular.4. Constructing and solving models: differential equations.5. Adjustment of parameters and bioprocesses optimization.6. Use of process simulators applied to biotechnology.6.1. Process design aided by simulation software6.2. Classification of process simulation software 6.3. Synthesis and analysis of process6.4. Process flowsheeting: concepts and limitation, convergence6.5. Application examples.'
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2).'
$ws.Range("C19").Value = 'Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2).'
$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'A média aritmética das notas corresponderá à média do período letivo, ou seja:Média do período letivo normal = (P1+ P2)/2.Serão aprovados os alunos que obtiverem média igual ou maior que 5,0.'
$ws.Range("C20").Value = 'A média aritmética das notas corresponderá à média do período letivo, ou seja:Média do período letivo normal = (P1+ P2)/2.Serão aprovados os alunos que obtiverem média igual ou maior que 5,0.'
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = 'Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0.'
$ws.Range("C21").Value = 'Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0.'
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'ADIDHARMA, H.; TEMYANKO, V. Mathcad for chemical engineers. Victoria, Canadá: Trafford Publishing, 2007. ISBN 1-4251-1541-1.BARRETO, L. S. Iniciação ao Scilab. 2 ed. Costa de Caparica, Portugal: Ebook, 2011. Disponível na internet: http://www.mat.ufrgs.br/~guidi/grad/MAT01169/SciLivro2.pdf Consulta em 17 de janeiro de 2014.CUTLIP, M. B.; SHACHAM, M. Problem solving in chemical and biochemical engineering with POLYMAT™, Excel, and MATLAB™. 2 ed. Boston, MA: Pearson Education, Inc., 2008. ISBN 978-0-13-148204-3.DIMIAN, A.C. Integrated design and simulation of chemical processes. Amsterdan, The Netherlands: Elsevier Science B.V., 2003, 678 p. ISBN 0-444-82996-2. ISSN: 1570-7946.ELNASHAIE, S. S. E. H.; GARHYAN, P. Conservation equations and modeling of chemical and biochemical processes. New York: Marcel dekker, Inc., 2003. ISBN 0-8247-0957-8.LOPES, L. C. O. Utilizando o SCILAB na Resolução de Problemas da Engenharia Química. v. 0.1. Curitiba/Parana,Brasil: XV COBEQ, Congresso Brasileiro de Engenharia Química, 2004. Disponível na internet: http://www.google.com.br/url?sa=t&rct=j&q=&esrc=s&source=web&cd=1&ved=0CC0QFjAA&url=http%3A%2F%2Fxa.yimg.com%2Fkq%2Fgroups%2F9656926%2F1161540061%2Fname%2Fscilab_COBEQ.pdf&ei=rnbZUpOZF4TJkAeHo4DwCg&usg=AFQjCNHc1deuW9_0qn7TyVvMEofHOUG9cA Consulta em 17 de janeiro de 2014.PINTO, J. C.; LAGE, P. L. C. Métodos numéricos em problemas de engenharia química. Rio de Janeiro, RJ: E-papers serviços Editoriais Ltda., 2001. ISBN 85-87922-11-4.RUGGIERO, M. A. G.; LOPES, V. L. R. Cálculo Numérico. Aspectos Teóricos e computacionais. 2 ed. São Paulo: Pearson education fo Brasil, 1998. ISBN 85-346-0204-2.SCHMIDELL, W.; LIMA, U. A.; AQUARONE, E.; BORZANI, W. Biotecnologia Industrial. Volume 2: Engenharia Bioquímica. São Paulo: Editora Edgard Blücher Ltda, 2001. ISBN 978-85-212-0279-0.SEIDER, W.D.; SEADER, J.D.; LEWIN, D.R.; WIDAGDO, S. Product and process design priciples. Synthesis, analysis, and Evaluation. 3 ed. Hoboken, NJ, USA: John Wiley & Sons, Inc., 2009, 728p. ISBN-13: 978-0470-04895-5.'
$ws.Range("C22").Value = 'ADIDHARMA, H.; TEMYANKO, V. Mathcad for chemical engineers. Victoria, Canadá: Trafford Publishing, 2007. ISBN 1-4251-1541-1.BARRETO, L. S. Iniciação ao Scilab. 2 ed. Costa de Caparica, Portugal: Ebook, 2011. Disponível na internet: http://www.mat.ufrgs.br/~guidi/grad/MAT01169/SciLivro2.pdf Consulta em 17 de janeiro de 2014.CUTLIP, M. B.; SHACHAM, M. Problem solving in chemical and biochemical engineering with POLYMAT™, Excel, and MATLAB™. 2 ed. Boston, MA: Pearson Education, Inc., 2008. ISBN 978-0-13-148204-3.DIMIAN, A.C. Integrated design and simulation of chemical processes. Amsterdan, The Netherlands: Elsevier Science B.V., 2003, 678 p. ISBN 0-444-82996-2. ISSN: 1570-7946.ELNASHAIE, S. S. E. H.; GARHYAN, P. Conservation equations and modeling of chemical and biochemical processes. New York: Marcel dekker, Inc., 2003. ISBN 0-8247-0957-8.LOPES, L. C. O. Utilizando o SCILAB na Resolução de Problemas da Engenharia Química. v. 0.1. Curitiba/Parana,Brasil: XV COBEQ, Congresso Brasileiro de Engenharia Química, 2004. Disponível na internet: http://www.google.com.br/url?sa=t&rct=j&q=&esrc=s&source=web&cd=1&ved=0CC0QFjAA&url=http%3A%2F%2Fxa.yimg.com%2Fkq%2Fgroups%2F9656926%2F1161540061%2Fname%2Fscilab_COBEQ.pdf&ei=rnbZUpOZF4TJkAeHo4DwCg&usg=AFQjCNHc1deuW9_0qn7TyVvMEofHOUG9cA Consulta em 17 de janeiro de 2014.PINTO, J. C.; LAGE, P. L. C. Métodos numéricos em problemas de engenharia química. Rio de Janeiro, RJ: E-papers serviços Editoriais Ltda., 2001. ISBN 85-87922-11-4.RUGGIERO, M. A. G.; LOPES, V. L. R. Cálculo Numérico. Aspectos Teóricos e computacionais. 2 ed. São Paulo: Pearson education fo Brasil, 1998. ISBN 85-346-0204-2.SCHMIDELL, W.; LIMA, U. A.; AQUARONE, E.; BORZANI, W. Biotecnologia Industrial. Volume 2: Engenharia Bioquímica. São Paulo: Editora Edgard Blücher Ltda, 2001. ISBN 978-85-212-0279-0.SEIDER, W.D.; SEADER, J.D.; LEWIN, D.R.; WIDAGDO, S. Product and process design priciples. Synthesis, analysis, and Evaluation. 3 ed. Hoboken, NJ, USA: John Wiley & Sons, Inc., 2009, 728p. ISBN-13: 978-0470-04895-5.'
$ws.Range("A23").Value = 'Requisitos:'
$ws.Range("B24").Value = 'LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)
'
$ws.Range("B25").Value = 'LOQ4086 -  Operações Unitárias II  (Requisito fraco)
'
$ws.Range("C25").Value = 'LOQ4086 -  Operações Unitárias II  (Requisito fraco)
'
